$wb = $excel.ActiveWorkbook

# Sheet "Person": A1 "name" -> "last_name"
$wsPerson = $wb.Worksheets.Item("Person")
$wsPerson.Range("A1").Value = "last_name"

# Sheet "Author": F1 "name" -> "last_name"
$wsAuthor = $wb.Worksheets.Item("Author")
$wsAuthor.Range("F1").Value = "last_name"

# Sheet "ImageSize": A1 "height" -> "height_im", B1 "width" -> "width_im"
$wsImageSize = $wb.Worksheets.Item("ImageSize")
$wsImageSize.Range("A1").Value = "height_im"
$wsImageSize.Range("B1").Value = "width_im"
